$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing data rows (2-29) before writing the regenerated simulation trace
$ws.Range("A2:I29").ClearContents()

# Row 2
$ws.Cells.Item(2,1).Value = "P1"
$ws.Cells.Item(2,2).Value = "J1"
$ws.Cells.Item(2,3).Value = "O11"
$ws.Cells.Item(2,4).Value = "M1"
$ws.Cells.Item(2,5).Value = "queued"
$ws.Cells.Item(2,6).Value = 0
$ws.Cells.Item(2,7).Value = 1
$ws.Cells.Item(2,8).Value = "O11"

# Row 3
$ws.Cells.Item(3,1).Value = "P1"
$ws.Cells.Item(3,2).Value = "J1"
$ws.Cells.Item(3,3).Value = "O11"
$ws.Cells.Item(3,4).Value = "M1"
$ws.Cells.Item(3,5).Value = "start"
$ws.Cells.Item(3,6).Value = 0
$ws.Cells.Item(3,7).Value = 0

# Row 4
$ws.Cells.Item(4,1).Value = "P2"
$ws.Cells.Item(4,2).Value = "J2"
$ws.Cells.Item(4,3).Value = "O21"
$ws.Cells.Item(4,4).Value = "M1"
$ws.Cells.Item(4,5).Value = "queued"
$ws.Cells.Item(4,6).Value = 2
$ws.Cells.Item(4,7).Value = 1
$ws.Cells.Item(4,8).Value = "O21"

# Row 5
$ws.Cells.Item(5,1).Value = "P1"
$ws.Cells.Item(5,2).Value = "J1"
$ws.Cells.Item(5,3).Value = "O11"
$ws.Cells.Item(5,4).Value = "M1"
$ws.Cells.Item(5,5).Value = "end"
$ws.Cells.Item(5,6).Value = 3.607996782156957

# Row 6
$ws.Cells.Item(6,1).Value = "P1"
$ws.Cells.Item(6,2).Value = "J1"
$ws.Cells.Item(6,3).Value = "O12"
$ws.Cells.Item(6,4).Value = "M1->M2"
$ws.Cells.Item(6,5).Value = "transfer"
$ws.Cells.Item(6,6).Value = 3.607996782156957
$ws.Cells.Item(6,9).Value = 1.668701124354598

# Row 7
$ws.Cells.Item(7,1).Value = "P2"
$ws.Cells.Item(7,2).Value = "J2"
$ws.Cells.Item(7,3).Value = "O21"
$ws.Cells.Item(7,4).Value = "M1"
$ws.Cells.Item(7,5).Value = "start"
$ws.Cells.Item(7,6).Value = 3.607996782156957
$ws.Cells.Item(7,7).Value = 0

# Row 8
$ws.Cells.Item(8,1).Value = "P3"
$ws.Cells.Item(8,2).Value = "J3"
$ws.Cells.Item(8,3).Value = "O31"
$ws.Cells.Item(8,4).Value = "M3"
$ws.Cells.Item(8,5).Value = "queued"
$ws.Cells.Item(8,6).Value = 4
$ws.Cells.Item(8,7).Value = 1
$ws.Cells.Item(8,8).Value = "O31"

# Row 9
$ws.Cells.Item(9,1).Value = "P3"
$ws.Cells.Item(9,2).Value = "J3"
$ws.Cells.Item(9,3).Value = "O31"
$ws.Cells.Item(9,4).Value = "M3"
$ws.Cells.Item(9,5).Value = "start"
$ws.Cells.Item(9,6).Value = 4
$ws.Cells.Item(9,7).Value = 0

# Row 10
$ws.Cells.Item(10,1).Value = "P1"
$ws.Cells.Item(10,2).Value = "J1"
$ws.Cells.Item(10,3).Value = "O12"
$ws.Cells.Item(10,4).Value = "M2"
$ws.Cells.Item(10,5).Value = "queued"
$ws.Cells.Item(10,6).Value = 5.276697906511554
$ws.Cells.Item(10,7).Value = 1
$ws.Cells.Item(10,8).Value = "O12"

# Row 11
$ws.Cells.Item(11,1).Value = "P1"
$ws.Cells.Item(11,2).Value = "J1"
$ws.Cells.Item(11,3).Value = "O12"
$ws.Cells.Item(11,4).Value = "M2"
$ws.Cells.Item(11,5).Value = "start"
$ws.Cells.Item(11,6).Value = 5.276697906511554
$ws.Cells.Item(11,7).Value = 0

# Row 12
$ws.Cells.Item(12,1).Value = "P2"
$ws.Cells.Item(12,2).Value = "J2"
$ws.Cells.Item(12,3).Value = "O21"
$ws.Cells.Item(12,4).Value = "M1"
$ws.Cells.Item(12,5).Value = "end"
$ws.Cells.Item(12,6).Value = 6.702257725566692

# Row 13
$ws.Cells.Item(13,1).Value = "P2"
$ws.Cells.Item(13,2).Value = "J2"
$ws.Cells.Item(13,3).Value = "O22"
$ws.Cells.Item(13,4).Value = "M1->M2"
$ws.Cells.Item(13,5).Value = "transfer"
$ws.Cells.Item(13,6).Value = 6.702257725566692
$ws.Cells.Item(13,9).Value = 1.7191080394563

# Row 14
$ws.Cells.Item(14,1).Value = "P2"
$ws.Cells.Item(14,2).Value = "J2"
$ws.Cells.Item(14,3).Value = "O22"
$ws.Cells.Item(14,4).Value = "M2"
$ws.Cells.Item(14,5).Value = "queued"
$ws.Cells.Item(14,6).Value = 8.421365765022992
$ws.Cells.Item(14,7).Value = 1
$ws.Cells.Item(14,8).Value = "O22"

# Row 15
$ws.Cells.Item(15,1).Value = "P3"
$ws.Cells.Item(15,2).Value = "J3"
$ws.Cells.Item(15,3).Value = "O31"
$ws.Cells.Item(15,4).Value = "M3"
$ws.Cells.Item(15,5).Value = "end"
$ws.Cells.Item(15,6).Value = 9.457171722217414

# Row 16
$ws.Cells.Item(16,1).Value = "P3"
$ws.Cells.Item(16,2).Value = "J3"
$ws.Cells.Item(16,3).Value = "O32"
$ws.Cells.Item(16,4).Value = "M3->M1"
$ws.Cells.Item(16,5).Value = "transfer"
$ws.Cells.Item(16,6).Value = 9.457171722217414
$ws.Cells.Item(16,9).Value = 1.690519551256694

# Row 17
$ws.Cells.Item(17,1).Value = "P1"
$ws.Cells.Item(17,2).Value = "J1"
$ws.Cells.Item(17,3).Value = "O12"
$ws.Cells.Item(17,4).Value = "M2"
$ws.Cells.Item(17,5).Value = "end"
$ws.Cells.Item(17,6).Value = 10.91628103788027

# Row 18
$ws.Cells.Item(18,1).Value = "P1"
$ws.Cells.Item(18,2).Value = "J1"
$ws.Cells.Item(18,3).Value = "O13"
$ws.Cells.Item(18,4).Value = "M2->M1"
$ws.Cells.Item(18,5).Value = "transfer"
$ws.Cells.Item(18,6).Value = 10.91628103788027
$ws.Cells.Item(18,9).Value = 0.719117402088626

# Row 19
$ws.Cells.Item(19,1).Value = "P2"
$ws.Cells.Item(19,2).Value = "J2"
$ws.Cells.Item(19,3).Value = "O22"
$ws.Cells.Item(19,4).Value = "M2"
$ws.Cells.Item(19,5).Value = "start"
$ws.Cells.Item(19,6).Value = 10.91628103788027
$ws.Cells.Item(19,7).Value = 0

# Row 20
$ws.Cells.Item(20,1).Value = "P3"
$ws.Cells.Item(20,2).Value = "J3"
$ws.Cells.Item(20,3).Value = "O32"
$ws.Cells.Item(20,4).Value = "M1"
$ws.Cells.Item(20,5).Value = "queued"
$ws.Cells.Item(20,6).Value = 11.14769127347411
$ws.Cells.Item(20,7).Value = 1
$ws.Cells.Item(20,8).Value = "O32"

# Row 21
$ws.Cells.Item(21,1).Value = "P3"
$ws.Cells.Item(21,2).Value = "J3"
$ws.Cells.Item(21,3).Value = "O32"
$ws.Cells.Item(21,4).Value = "M1"
$ws.Cells.Item(21,5).Value = "start"
$ws.Cells.Item(21,6).Value = 11.14769127347411
$ws.Cells.Item(21,7).Value = 0

# Row 22
$ws.Cells.Item(22,1).Value = "P1"
$ws.Cells.Item(22,2).Value = "J1"
$ws.Cells.Item(22,3).Value = "O13"
$ws.Cells.Item(22,4).Value = "M1"
$ws.Cells.Item(22,5).Value = "queued"
$ws.Cells.Item(22,6).Value = 11.63539843996889
$ws.Cells.Item(22,7).Value = 1
$ws.Cells.Item(22,8).Value = "O13"

# Row 23
$ws.Cells.Item(23,1).Value = "P3"
$ws.Cells.Item(23,2).Value = "J3"
$ws.Cells.Item(23,3).Value = "O32"
$ws.Cells.Item(23,4).Value = "M1"
$ws.Cells.Item(23,5).Value = "end"
$ws.Cells.Item(23,6).Value = 14.54747669087936

# Row 24
$ws.Cells.Item(24,1).Value = "P3"
$ws.Cells.Item(24,2).Value = "J3"
$ws.Cells.Item(24,5).Value = "done"
$ws.Cells.Item(24,6).Value = 14.54747669087936

# Row 25
$ws.Cells.Item(25,1).Value = "P1"
$ws.Cells.Item(25,2).Value = "J1"
$ws.Cells.Item(25,3).Value = "O13"
$ws.Cells.Item(25,4).Value = "M1"
$ws.Cells.Item(25,5).Value = "start"
$ws.Cells.Item(25,6).Value = 14.54747669087936
$ws.Cells.Item(25,7).Value = 0

# Row 26
$ws.Cells.Item(26,1).Value = "P2"
$ws.Cells.Item(26,2).Value = "J2"
$ws.Cells.Item(26,3).Value = "O22"
$ws.Cells.Item(26,4).Value = "M2"
$ws.Cells.Item(26,5).Value = "end"
$ws.Cells.Item(26,6).Value = 14.58431033061374

# Row 27
$ws.Cells.Item(27,1).Value = "P2"
$ws.Cells.Item(27,2).Value = "J2"
$ws.Cells.Item(27,5).Value = "done"
$ws.Cells.Item(27,6).Value = 14.58431033061374

# Row 28
$ws.Cells.Item(28,1).Value = "P1"
$ws.Cells.Item(28,2).Value = "J1"
$ws.Cells.Item(28,3).Value = "O13"
$ws.Cells.Item(28,4).Value = "M1"
$ws.Cells.Item(28,5).Value = "end"
$ws.Cells.Item(28,6).Value = 17.05990858031197

# Row 29
$ws.Cells.Item(29,1).Value = "P1"
$ws.Cells.Item(29,2).Value = "J1"
$ws.Cells.Item(29,5).Value = "done"
$ws.Cells.Item(29,6).Value = 17.05990858031197
